$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row continues the table directly below the last data row (row 4),
# so copy its formatting down to row 5 first (keeps the same cell style
# used by the other data rows: vertical-centered, wrapped text).
$ws.Range("A4:C4").Copy()
$ws.Range("A5:C5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the new record: ID=4, Gender=M, Ht=188 cm
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "M"
$ws.Range("C5").Value = "188 cm"

# Leave the selection where Excel would land after typing the row and
# pressing Enter - the first cell of the next empty row.
$ws.Range("A6").Select()
